# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# cryptocurrency tracking sheet. Values that look numeric are written with a
# leading apostrophe so Excel keeps them as text (matching the source data,
# which stores prices/percentages as formatted strings, e.g. "10.40" must
# stay "10.40" and not become the number 10.4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.813.88"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.266.52"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'303.74"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'92.65"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'32.42"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D12").Value = "'0.113"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "2.615.88"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "2.265.56"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  +3.93%  "
$ws.Range("D18").Value = "41.755.09"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'12.82"
$ws.Range("E19").Value = "  +4.11%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'5.94"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'67.64"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "'244.22"
$ws.Range("D24").Value = "'2.59"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'1.93"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "'24.01"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "'9.58"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -5.75%  "
$ws.Range("D30").Value = "'34.96"
$ws.Range("D31").Value = "'159.06"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'3.02"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'16.89"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "'3.94"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").Value = "'19.96"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").Value = "2.008.28"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "'2.27"
$ws.Range("E44").Value = "  +13.22%  "
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").Value = "'10.40"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D47").Value = "'2.90"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'53.26"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "'73.13"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  +0.16%  "
